$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    ,("B2", 7.82133082566617)
    ,("D2", 9.988398897164046)
    ,("E2", 14.32816096069627)
    ,("F2", 29.90354716093452)
    ,("G2", 29.17753511029934)
    ,("H2", 14.49866275347064)
    ,("J2", 10.21976011310801)
    ,("K2", 10.29665410415161)
    ,("M2", 14.57991496661228)
    ,("O2", 22.10560947689195)
    ,("B3", 7.750599921545449)
    ,("D3", 9.963905134536027)
    ,("E3", 14.34425216291507)
    ,("F3", 29.99830103918405)
    ,("G3", 29.30188088489006)
    ,("H3", 14.55387829129145)
    ,("J3", 10.2482626535018)
    ,("K3", 9.849906929832551)
    ,("M3", 14.41233242176732)
    ,("O3", 22.20127598429987)
    ,("B4", 7.708686723291977)
    ,("D4", 9.950364002103923)
    ,("E4", 14.35666734359589)
    ,("F4", 30.06388129717048)
    ,("G4", 29.3880641384802)
    ,("H4", 14.59016273145722)
    ,("J4", 10.26701508742808)
    ,("K4", 9.563577937218119)
    ,("M4", 14.30980147068305)
    ,("O4", 22.26492279276048)
    ,("B5", 7.692006826993126)
    ,("D5", 9.945226605794771)
    ,("E5", 14.36236434824852)
    ,("F5", 30.09246230383482)
    ,("G5", 29.42564550237796)
    ,("H5", 14.6055480778826)
    ,("J5", 10.27497198996644)
    ,("K5", 9.443984391384275)
    ,("M5", 14.26814992408395)
    ,("O5", 22.2920914241658)
    ,("B6", 7.689261843937516)
    ,("D6", 9.944396651066482)
    ,("E6", 14.3633488520262)
    ,("F6", 30.09732014974185)
    ,("G6", 29.43203413792591)
    ,("H6", 14.60813899580017)
    ,("J6", 10.2763122717398)
    ,("K6", 9.423953673097465)
    ,("M6", 14.26124273670842)
    ,("O6", 22.29667711093184)
    ,("B7", 7.70846012807757)
    ,("D7", 9.950293170764525)
    ,("E7", 14.35674159309273)
    ,("F7", 30.0642592398352)
    ,("G7", 29.38856102562008)
    ,("H7", 14.59036779744346)
    ,("J7", 10.26712112054269)
    ,("K7", 9.561976688552848)
    ,("M7", 14.30923916329059)
    ,("O7", 22.26528421247818)
    ,("B8", 7.796639748558013)
    ,("D8", 9.979645000428308)
    ,("E8", 14.33318323371869)
    ,("F8", 29.93467958767732)
    ,("G8", 29.21836134654646)
    ,("H8", 14.51720683239567)
    ,("J8", 10.22932820836047)
    ,("K8", 10.14517321975912)
    ,("M8", 14.52207995259046)
    ,("O8", 22.13757565847529)
    ,("B9", 7.980724238261707)
    ,("D9", 10.04890274879951)
    ,("E9", 14.30708455521427)
    ,("F9", 29.7394887712873)
    ,("G9", 28.963144773405)
    ,("H9", 14.3926284021754)
    ,("J9", 10.16513199650206)
    ,("K9", 11.18930472606495)
    ,("M9", 14.94065659480468)
    ,("O9", 21.92616068597454)
    ,("B10", 8.121583581056548)
    ,("D10", 10.10664903526043)
    ,("E10", 14.30013372360042)
    ,("F10", 29.63224183763653)
    ,("G10", 28.82418602670313)
    ,("H10", 14.31260374304968)
    ,("J10", 10.12398899478217)
    ,("K10", 11.89123028698282)
    ,("M10", 15.24664801105621)
    ,("O10", 21.79473741221788)
    ,("B11", 8.186621758304877)
    ,("D11", 10.134346301941)
    ,("E11", 14.29961616164746)
    ,("F11", 29.59134962830333)
    ,("G11", 28.77164208254204)
    ,("H11", 14.27869383829717)
    ,("J11", 10.10657480685216)
    ,("K11", 12.19567679407965)
    ,("M11", 15.38504528143899)
    ,("O11", 21.74016447549718)
    ,("B12", 8.211366059068775)
    ,("D12", 10.14503412381946)
    ,("E12", 14.29979926929119)
    ,("F12", 29.5770032325083)
    ,("G12", 28.75328889437468)
    ,("H12", 14.26621153114499)
    ,("J12", 10.10016737837107)
    ,("K12", 12.30877810385911)
    ,("M12", 15.43729906783768)
    ,("O12", 21.72025070949988)
    ,("B13", 8.206032138651119)
    ,("D13", 10.14272353539608)
    ,("E13", 14.29974299394936)
    ,("F13", 29.5800422977039)
    ,("G13", 28.75717276022998)
    ,("H13", 14.26888386795036)
    ,("J13", 10.10153902258053)
    ,("K13", 12.28451762299746)
    ,("M13", 15.42605275531546)
    ,("O13", 21.72450602121621)
    ,("B14", 8.188655294911943)
    ,("D14", 10.13522163254099)
    ,("E14", 14.2996236357127)
    ,("F14", 29.59014650497928)
    ,("G14", 28.77010115557895)
    ,("H14", 14.27765972343275)
    ,("J14", 10.10604391872937)
    ,("K14", 12.20502576669777)
    ,("M14", 15.38934754483882)
    ,("O14", 21.73851107718097)
    ,("B15", 8.178025887735178)
    ,("D15", 10.13065229759033)
    ,("E15", 14.29959985710531)
    ,("F15", 29.59648399398647)
    ,("G15", 28.77822153345391)
    ,("H15", 14.2830818951153)
    ,("J15", 10.10882763555732)
    ,("K15", 12.15604872729816)
    ,("M15", 15.36684330619087)
    ,("O15", 21.74718755600437)
    ,("B16", 8.117350599126754)
    ,("D16", 10.10486719365243)
    ,("E16", 14.3002206644983)
    ,("F16", 29.63507339855157)
    ,("G16", 28.82783544635948)
    ,("H16", 14.31487001848241)
    ,("J16", 10.12515322827168)
    ,("K16", 11.87103055528228)
    ,("M16", 15.23758395411025)
    ,("O16", 21.79840896270292)
    ,("B17", 8.080358376295008)
    ,("D17", 10.0894105413516)
    ,("E17", 14.30127803766098)
    ,("F17", 29.66077131810787)
    ,("G17", 28.86101183586284)
    ,("H17", 14.33500972993278)
    ,("J17", 10.13550171932559)
    ,("K17", 11.69233700977404)
    ,("M17", 15.1580544727142)
    ,("O17", 21.83116825346947)
    ,("B18", 8.059173181903336)
    ,("D18", 10.08065508166781)
    ,("E18", 14.30213519667325)
    ,("F18", 29.67629501928254)
    ,("G18", 28.88109769856661)
    ,("H18", 14.34682827585512)
    ,("J18", 10.14157646143995)
    ,("K18", 11.58816056888051)
    ,("M18", 15.1122382354537)
    ,("O18", 21.85050098775537)
    ,("B19", 8.05201667978708)
    ,("D19", 10.07771396088036)
    ,("E19", 14.30246821190708)
    ,("F19", 29.68167857064021)
    ,("G19", 28.88807049709584)
    ,("H19", 14.35087014410807)
    ,("J19", 10.14365432363639)
    ,("K19", 11.55264996499237)
    ,("M19", 15.0967143209602)
    ,("O19", 21.85713088489061)
    ,("B20", 8.084286927542468)
    ,("D20", 10.0910420192749)
    ,("E20", 14.30113971720622)
    ,("F20", 29.65795881419975)
    ,("G20", 28.85737620153176)
    ,("H20", 14.33284153069761)
    ,("J20", 10.13438742190195)
    ,("K20", 11.71150416905149)
    ,("M20", 15.16652836254887)
    ,("O20", 21.82763019375571)
    ,("B21", 8.193756331665737)
    ,("D21", 10.13741975803469)
    ,("E21", 14.29964841577741)
    ,("F21", 29.58714772879403)
    ,("G21", 28.76626179219696)
    ,("H21", 14.27507230819282)
    ,("J21", 10.10471565002794)
    ,("K21", 12.22843412570641)
    ,("M21", 15.40013324943285)
    ,("O21", 21.73437703019796)
    ,("B22", 8.265965301905466)
    ,("D22", 10.16889003370698)
    ,("E22", 14.30088275506473)
    ,("F22", 29.54750645198952)
    ,("G22", 28.71571727303317)
    ,("H22", 14.23940722825986)
    ,("J22", 10.08641297448832)
    ,("K22", 12.55351921336291)
    ,("M22", 15.55189092262548)
    ,("O22", 21.6778139428362)
    ,("B23", 8.227372597688811)
    ,("D23", 10.15198965261691)
    ,("E23", 14.30002227704925)
    ,("F23", 29.56805539319046)
    ,("G23", 28.74186690533774)
    ,("H23", 14.25825105820982)
    ,("J23", 10.09608185354925)
    ,("K23", 12.38119659277155)
    ,("M23", 15.47099158061542)
    ,("O23", 21.70760090000057)
    ,("B24", 8.082510572010507)
    ,("D24", 10.09030402031868)
    ,("E24", 14.30120147538266)
    ,("F24", 29.65922801244539)
    ,("G24", 28.85901671928078)
    ,("H24", 14.33382102611623)
    ,("J24", 10.13489080558224)
    ,("K24", 11.70284318766219)
    ,("M24", 15.16269760574493)
    ,("O24", 21.82922819642524)
    ,("B25", 7.929850910170794)
    ,("D25", 10.02894089964039)
    ,("E25", 14.31199547957542)
    ,("F25", 29.78595934879031)
    ,("G25", 29.02370778239584)
    ,("H25", 14.42430923712234)
    ,("J25", 10.18143962000373)
    ,("K25", 10.91800932796041)
    ,("M25", 14.82752780207331)
    ,("O25", 21.97916464603113)
)

foreach ($item in $changes) {
    $addr = $item[0]
    $val = $item[1]
    $ws.Range($addr).Value = $val
}